$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198, pushing existing rows 198:285 down to 199:286
$ws.Rows.Item(198).Insert()

# Populate the new row 198 with values
$ws.Cells.Item(198, 1).Value = 4
$ws.Cells.Item(198, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(198, 3).Value = "Los Lagos"
$ws.Cells.Item(198, 4).Value = 44755
$ws.Cells.Item(198, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(198, 5).Value = 10
$ws.Cells.Item(198, 6).Value = "Fruta"
$ws.Cells.Item(198, 7).Value = 100104
$ws.Cells.Item(198, 8).Value = "Frutos de pepita"
$ws.Cells.Item(198, 9).Value = 100104005
$ws.Cells.Item(198, 10).Value = "Pera"
$ws.Cells.Item(198, 11).Value = "Packham's Triumph"
$ws.Cells.Item(198, 12).Value = "Primera"
$ws.Cells.Item(198, 13).Value = 200
$ws.Cells.Item(198, 14).Value = 16000
$ws.Cells.Item(198, 15).Value = 16000
$ws.Cells.Item(198, 16).Value = 16000
$ws.Cells.Item(198, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(198, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(198, 19).Value = 1067
$ws.Cells.Item(198, 20).Value = 15
